$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "hi"
$ws.Range("B1").Value = "hello"
$ws.Range("A2").Value = "hello"
$ws.Range("B2").Value = "hi"

$ws.Range("B2").Select()
